$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 114.949946438405
$ws.Range("B1").Value = 185.887906973198
$ws.Range("C1").Value = 125.428022316391
$ws.Range("D1").Value = 37.5679961580634
$ws.Range("E1").Value = 185.454230096868
$ws.Range("F1").Value = 140.754750622788
$ws.Range("G1").Value = 123.34300536818
$ws.Range("H1").Value = 100.982537540133
$ws.Range("I1").Value = 13.9449718473223
$ws.Range("J1").Value = 101.5689083848
$ws.Range("A2").Value = 108.717175530604
$ws.Range("B2").Value = 67.9643290433867
$ws.Range("C2").Value = 84.8974579409219
$ws.Range("D2").Value = 92.8327521741543
$ws.Range("E2").Value = 17.196502917072
$ws.Range("F2").Value = 197.092276344584
$ws.Range("G2").Value = 77.6451101888181
$ws.Range("H2").Value = 107.689645377775
$ws.Range("I2").Value = 178.059247870957
$ws.Range("J2").Value = 34.7992823621255
$ws.Range("A3").Value = 60.0700568687497
$ws.Range("B3").Value = 20.3738832941157
$ws.Range("C3").Value = 111.097120824781
$ws.Range("D3").Value = 126.994203555861
$ws.Range("E3").Value = 125.388765114075
$ws.Range("F3").Value = 115.580664535789
$ws.Range("G3").Value = 186.699653969472
$ws.Range("H3").Value = 151.489565219492
$ws.Range("I3").Value = 179.992093881588
$ws.Range("J3").Value = 153.088201188058
$ws.Range("A4").Value = 155.470422262079
$ws.Range("B4").Value = 90.3641691852194
$ws.Range("C4").Value = 182.954371898879
$ws.Range("D4").Value = 192.718372397459
$ws.Range("E4").Value = 55.0380260939887
$ws.Range("F4").Value = 58.4344572659742
$ws.Range("G4").Value = 185.548592538363
$ws.Range("H4").Value = 85.797343629318
$ws.Range("I4").Value = 81.7145861134467
$ws.Range("J4").Value = 138.018755865292
$ws.Range("A5").Value = 152.584969975327
$ws.Range("B5").Value = 82.5259498704811
$ws.Range("C5").Value = 94.7762802684569
$ws.Range("D5").Value = 54.5953426764325
$ws.Range("E5").Value = 36.9117616847678
$ws.Range("F5").Value = 87.587959918933
$ws.Range("G5").Value = 171.969674700857
$ws.Range("H5").Value = 102.146972949685
$ws.Range("I5").Value = 97.8994711758101
$ws.Range("J5").Value = 139.735567355405
$ws.Range("A6").Value = 133.168626173012
$ws.Range("B6").Value = 163.818816544404
$ws.Range("C6").Value = 124.715224245896
$ws.Range("D6").Value = 16.2769996636906
$ws.Range("E6").Value = 109.765918976518
$ws.Range("F6").Value = 94.576063144289
$ws.Range("G6").Value = 74.7907861484172
$ws.Range("H6").Value = 198.43381876053
$ws.Range("I6").Value = 112.179231043988
$ws.Range("J6").Value = 69.8735655610792
$ws.Range("A7").Value = 154.055096653316
$ws.Range("B7").Value = 171.853440148688
$ws.Range("C7").Value = 120.990443658545
$ws.Range("D7").Value = 60.8567706592645
$ws.Range("E7").Value = 146.098486122721
$ws.Range("F7").Value = 18.3530063453843
$ws.Range("G7").Value = 85.0099571445072
$ws.Range("H7").Value = 92.1790855434626
$ws.Range("I7").Value = 37.7947260801656
$ws.Range("J7").Value = 158.762045651098
$ws.Range("A8").Value = 11.5436838062218
$ws.Range("B8").Value = 191.8477665595
$ws.Range("C8").Value = 25.9750592643279
$ws.Range("D8").Value = 40.0404920056651
$ws.Range("E8").Value = 82.2143123867988
$ws.Range("F8").Value = 177.544106998269
$ws.Range("G8").Value = 125.597603025659
$ws.Range("H8").Value = 56.5017781483483
$ws.Range("I8").Value = 90.0824418710929
$ws.Range("J8").Value = 37.800805195142
$ws.Range("A9").Value = 143.610989834932
$ws.Range("B9").Value = 84.5526810197871
$ws.Range("C9").Value = 53.5900940436824
$ws.Range("D9").Value = 40.2565265261831
$ws.Range("E9").Value = 19.9195750150455
$ws.Range("F9").Value = 191.651605717676
$ws.Range("G9").Value = 165.648944939323
$ws.Range("H9").Value = 166.677372235189
$ws.Range("I9").Value = 82.952453420941
$ws.Range("J9").Value = 160.4619629497
$ws.Range("A10").Value = 183.643671117557
$ws.Range("B10").Value = 187.114773777833
$ws.Range("C10").Value = 173.618112585329
$ws.Range("D10").Value = 11.8410205523674
$ws.Range("E10").Value = 183.963659211976
$ws.Range("F10").Value = 180.731529826639
$ws.Range("G10").Value = 161.535506211936
$ws.Range("H10").Value = 33.9195096091924
$ws.Range("I10").Value = 108.496856553711
$ws.Range("J10").Value = 18.5587553393835
$ws.Range("A11").Value = 2.57800277442578
$ws.Range("B11").Value = 79.7905891573944
$ws.Range("C11").Value = 64.3522468695195
$ws.Range("D11").Value = 139.137425524712
$ws.Range("E11").Value = 128.191883549184
$ws.Range("F11").Value = 141.320859613512
$ws.Range("G11").Value = 137.843757280076
$ws.Range("H11").Value = 84.6747322402311
$ws.Range("I11").Value = 134.062687276892
$ws.Range("J11").Value = 132.22181197825
$ws.Range("A12").Value = 168.97846011863
$ws.Range("B12").Value = 18.2890080000688
$ws.Range("C12").Value = 108.351376889437
$ws.Range("D12").Value = 74.3784258488465
$ws.Range("E12").Value = 126.262575726147
$ws.Range("F12").Value = 69.5024156335287
$ws.Range("G12").Value = 118.263346105005
$ws.Range("H12").Value = 80.7339171323618
$ws.Range("I12").Value = 40.937195644219
$ws.Range("J12").Value = 154.446880405046
$ws.Range("A13").Value = 52.7040614060611
$ws.Range("B13").Value = 118.332584909318
$ws.Range("C13").Value = 9.22663212252158
$ws.Range("D13").Value = 77.3327631304659
$ws.Range("E13").Value = 175.118374533541
$ws.Range("F13").Value = 24.4289100283891
$ws.Range("G13").Value = 18.2296539741706
$ws.Range("H13").Value = 14.1340387119604
$ws.Range("I13").Value = 56.0768327936888
$ws.Range("J13").Value = 101.48278256016
$ws.Range("A14").Value = 16.0086007863323
$ws.Range("B14").Value = 91.6780934164664
$ws.Range("C14").Value = 148.004921594637
$ws.Range("D14").Value = 71.5236865317094
$ws.Range("E14").Value = 35.2228024207162
$ws.Range("F14").Value = 63.8204006775377
$ws.Range("G14").Value = 20.2004341502676
$ws.Range("H14").Value = 114.45266851897
$ws.Range("I14").Value = 112.064642977
$ws.Range("J14").Value = 78.5987154015334
$ws.Range("A15").Value = 53.8078484375998
$ws.Range("B15").Value = 80.9742126990921
$ws.Range("C15").Value = 32.6146849582971
$ws.Range("D15").Value = 150.730641442691
$ws.Range("E15").Value = 191.48350283107
$ws.Range("F15").Value = 165.354663117488
$ws.Range("G15").Value = 78.7633968883955
$ws.Range("H15").Value = 99.239686736483
$ws.Range("I15").Value = 85.5784448262204
$ws.Range("J15").Value = 114.461243578448
$ws.Range("A16").Value = 62.4681837216337
$ws.Range("B16").Value = 80.8015890795745
$ws.Range("C16").Value = 192.982313964973
$ws.Range("D16").Value = 154.049976148666
$ws.Range("E16").Value = 165.854693933322
$ws.Range("F16").Value = 84.2454178651075
$ws.Range("G16").Value = 70.5639570348728
$ws.Range("H16").Value = 187.019483739054
$ws.Range("I16").Value = 164.019050991172
$ws.Range("J16").Value = 103.762973520794
$ws.Range("A17").Value = 123.091205639341
$ws.Range("B17").Value = 123.709718568115
$ws.Range("C17").Value = 28.5978994465423
$ws.Range("D17").Value = 32.5799047167319
$ws.Range("E17").Value = 116.213211191917
$ws.Range("F17").Value = 77.3003667021638
$ws.Range("G17").Value = 70.2840864054319
$ws.Range("H17").Value = 36.8276903577278
$ws.Range("I17").Value = 39.1556234281303
$ws.Range("J17").Value = 62.4421750486094
$ws.Range("A18").Value = 49.3019814832611
$ws.Range("B18").Value = 3.81067758603519
$ws.Range("C18").Value = 168.669274155362
$ws.Range("D18").Value = 162.338480242686
$ws.Range("E18").Value = 100.639031967446
$ws.Range("F18").Value = 171.729848706969
$ws.Range("G18").Value = 85.7178999510211
$ws.Range("H18").Value = 58.4959906798303
$ws.Range("I18").Value = 85.8492602993964
$ws.Range("J18").Value = 9.76371141605252
$ws.Range("A19").Value = 145.665513139994
$ws.Range("B19").Value = 118.989967237688
$ws.Range("C19").Value = 128.55559388574
$ws.Range("D19").Value = 141.615589215241
$ws.Range("E19").Value = 39.0145988385261
$ws.Range("F19").Value = 135.207011706758
$ws.Range("G19").Value = 98.6957794514931
$ws.Range("H19").Value = 193.954945445971
$ws.Range("I19").Value = 105.668992598387
$ws.Range("J19").Value = 150.977384555609
$ws.Range("A20").Value = 193.256443642665
$ws.Range("B20").Value = 33.180950411214
$ws.Range("C20").Value = 150.433617527799
$ws.Range("D20").Value = 8.30166945620518
$ws.Range("E20").Value = 155.507509762192
$ws.Range("F20").Value = 130.098129869484
$ws.Range("G20").Value = 52.3763132525498
$ws.Range("H20").Value = 0.0347802415652109
$ws.Range("I20").Value = 34.5174302507739
$ws.Range("J20").Value = 114.183136128906
